# Fruta / hortaliza, semanal
# Insert two new weekly price observations (rows 93-94) for
# "Feria Lagunitas de Puerto Montt - Mandarina", pushing the existing
# historical rows (old 93-108) down to (95-110).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows at row 93, shifting everything from the old row 93
# onward down by two rows (Excel default Insert shifts cells down).
$ws.Rows("93:94").Insert()

# --- New row 93: Mandarina / Murcott / Primera ---
$ws.Cells.Item(93, 1).Value = 4
$ws.Cells.Item(93, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(93, 3).Value = "Los Lagos"
$ws.Cells.Item(93, 4).Value = 44476
$ws.Cells.Item(93, 5).Value = 10
$ws.Cells.Item(93, 6).Value = "Fruta"
$ws.Cells.Item(93, 7).Value = 100102
$ws.Cells.Item(93, 8).Value = "Cítricos"
$ws.Cells.Item(93, 9).Value = 100102004
$ws.Cells.Item(93, 10).Value = "Mandarina"
$ws.Cells.Item(93, 11).Value = "Murcott"
$ws.Cells.Item(93, 12).Value = "Primera"
$ws.Cells.Item(93, 13).Value = 600
$ws.Cells.Item(93, 14).Value = 6000
$ws.Cells.Item(93, 15).Value = 6500
$ws.Cells.Item(93, 16).Value = 6250
$ws.Cells.Item(93, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(93, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(93, 19).Value = 625
$ws.Cells.Item(93, 20).Value = 10

# --- New row 94: Mandarina / Murcott / Segunda ---
$ws.Cells.Item(94, 1).Value = 4
$ws.Cells.Item(94, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(94, 3).Value = "Los Lagos"
$ws.Cells.Item(94, 4).Value = 44476
$ws.Cells.Item(94, 5).Value = 10
$ws.Cells.Item(94, 6).Value = "Fruta"
$ws.Cells.Item(94, 7).Value = 100102
$ws.Cells.Item(94, 8).Value = "Cítricos"
$ws.Cells.Item(94, 9).Value = 100102004
$ws.Cells.Item(94, 10).Value = "Mandarina"
$ws.Cells.Item(94, 11).Value = "Murcott"
$ws.Cells.Item(94, 12).Value = "Segunda"
$ws.Cells.Item(94, 13).Value = 200
$ws.Cells.Item(94, 14).Value = 4500
$ws.Cells.Item(94, 15).Value = 4500
$ws.Cells.Item(94, 16).Value = 4500
$ws.Cells.Item(94, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(94, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(94, 19).Value = 450
$ws.Cells.Item(94, 20).Value = 10
